# Metadata cleanup: fill in the "sample" rows that were mistakenly left as
# blank/NA placeholders at the edges of each blank-well block, and correct
# the block boundary so the blank-well range lines up with the physical
# plate layout (size-normalized analysis).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# plate1 / 16C / control block (rows 42-56): first blank row (43) should
# match the "control-1 / right" header row above it (42); last blank row
# (55) should match the "control-2 / left" header row below it (56).
$ws.Range("D43").Value = "control"
$ws.Range("E43").Value = "control-1"
$ws.Range("F43").Value = "right"
$ws.Range("H43").Value = "sample"

$ws.Range("D55").Value = "control"
$ws.Range("E55").Value = "control-2"
$ws.Range("F55").Value = "left"
$ws.Range("H55").Value = "sample"

# plate2 / 16C / treated block (rows 138-152): shift the blank range down
# by one row - row 139 becomes a "treated-4 / left" sample row, and row
# 151 becomes a blank/NA row instead of "treated-5 / right".
$ws.Range("E139").Value = "treated-4"
$ws.Range("F139").Value = "left"
$ws.Range("H139").Value = "sample"

$ws.Range("D140").Value = "NA"
$ws.Range("D141").Value = "NA"
$ws.Range("D142").Value = "NA"
$ws.Range("D143").Value = "NA"
$ws.Range("D144").Value = "NA"
$ws.Range("D145").Value = "NA"
$ws.Range("D146").Value = "NA"
$ws.Range("D147").Value = "NA"
$ws.Range("D148").Value = "NA"
$ws.Range("D149").Value = "NA"
$ws.Range("D150").Value = "NA"
$ws.Range("D151").Value = "NA"
$ws.Range("E151").Value = "NA"
$ws.Range("F151").Value = "NA"
$ws.Range("H151").Value = "blank"

# plate3 / 40C / control block (rows 234-248): first blank row (235)
# should match the "control-1 / right" header row above it (234).
$ws.Range("D235").Value = "control"
$ws.Range("E235").Value = "control-1"
$ws.Range("F235").Value = "right"
$ws.Range("H235").Value = "sample"

# plate4 / 40C / treated block (rows 330-344): shift the blank range down
# by one row - row 331 becomes a "treated-4 / left" sample row, and rows
# 332-343 become fully blank/NA (spat_treatment column included).
$ws.Range("E331").Value = "treated-4"
$ws.Range("F331").Value = "left"
$ws.Range("H331").Value = "sample"

$ws.Range("D332").Value = "NA"
$ws.Range("D333").Value = "NA"
$ws.Range("D334").Value = "NA"
$ws.Range("D335").Value = "NA"
$ws.Range("D336").Value = "NA"
$ws.Range("D337").Value = "NA"
$ws.Range("D338").Value = "NA"
$ws.Range("D339").Value = "NA"
$ws.Range("D340").Value = "NA"
$ws.Range("D341").Value = "NA"
$ws.Range("D342").Value = "NA"
$ws.Range("D343").Value = "NA"

# Restore the view/selection state recorded when the file was last saved.
$ws.Range("F149").Select()
